# Applies the member-list update captured by the diff:
#  - Row 76: rename the member from "Mame Faty Laye / THIAW" (with an address
#    of "CROISEMENT MALIKA") to "Ousmane / SARR" and clear the now-unused
#    address cell.
#  - Row 77: rename the member from "Astou Laye / Sow" to "Aiisata / LY".
#  - Rows 78-80: three new members are appended - "Assane SALL",
#    "MAMADOU LAMINE  SANE" and "Daouda THIOUNE".
#  - Misc view-state bookkeeping (active cell / scroll position) is updated
#    to match where the author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste des membres")

# Row 76 - existing member renamed, address cleared
$ws.Range("A76").Value = "Ousmane "
$ws.Range("B76").Value = "SARR"
$ws.Range("C76").ClearContents()

# Row 77 - existing member renamed
$ws.Range("A77").Value = "Aiisata "
$ws.Range("B77").Value = "LY"

# Row 78 - new member
$ws.Range("A78").Value = "Assane "
$ws.Range("B78").Value = "SALL"

# Row 79 - new member
$ws.Range("A79").Value = "MAMADOU LAMINE "
$ws.Range("B79").Value = "SANE"

# Row 80 - new member
$ws.Range("A80").Value = "Daouda"
$ws.Range("B80").Value = "THIOUNE"

# Restore the author's final selection / scroll position
$null = $ws.Range("H81").Select()
